$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 102, shifting existing rows 102:210 down to 103:210
$ws.Rows(102).Insert()

# Populate the newly inserted row 102 with the new record
$ws.Cells.Item(102, 1).Value = 10
$ws.Cells.Item(102, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value = "La Araucanía"
$ws.Cells.Item(102, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(102, 5).Value = 9
$ws.Cells.Item(102, 6).Value = 100112009
$ws.Cells.Item(102, 7).Value = "Acelga"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 40
$ws.Cells.Item(102, 11).Value = 7000
$ws.Cells.Item(102, 12).Value = 7000
$ws.Cells.Item(102, 13).Value = 7000
$ws.Cells.Item(102, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(102, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(102, 16).Value = 583
$ws.Cells.Item(102, 17).Value = 12
$ws.Cells.Item(102, 18).Value = "Hortaliza"
